$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.218.07"
$ws.Range("E2").Value = "  -0.98%  "
$ws.Range("D3").Value = "1.865.86"
$ws.Range("E3").Value = "  -0.71%  "
$ws.Range("D4").Value = "0.9995"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7100"
$ws.Range("E5").Value = "  -0.71%  "
$ws.Range("D6").Value = "241.74"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").Value = "0.3114"
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").Value = "0.07662"
$ws.Range("E9").Value = "  -3.63%  "
$ws.Range("D10").Value = "24.68"
$ws.Range("E10").Value = "  -2.41%  "
$ws.Range("D11").Value = "0.08369"
$ws.Range("E11").Value = "  +0.96%  "
$ws.Range("D12").Value = "1.859.26"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").Value = "5.224"
$ws.Range("E13").Value = "  -1.19%  "
$ws.Range("D14").Value = "0.7108"
$ws.Range("E14").Value = "  -2.61%  "
$ws.Range("D15").Value = "91.37"
$ws.Range("E15").Value = "  +0.10%  "
$ws.Range("D16").Value = "29.232.24"
$ws.Range("D17").Value = "5.947"
$ws.Range("E17").Value = "  +0.19%  "
$ws.Range("D18").Value = "243.81"
$ws.Range("E18").Value = "  -1.14%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007832"
$ws.Range("E19").Value = "  -0.72%  "
$ws.Range("D20").Value = "2.114.78"
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("D21").Value = "13.09"
$ws.Range("E21").Value = "  -1.98%  "
$ws.Range("D22").Value = "0.9994"
$ws.Range("D23").Value = "7.867"
$ws.Range("E23").Value = "  -1.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.0000"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1650"
$ws.Range("E25").Value = "  +1.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.00"
$ws.Range("E26").Value = "  -0.55%  "
$ws.Range("D27").Value = "8.956"
$ws.Range("E27").Value = "  -1.17%  "
$ws.Range("E28").Value = "  +0.83%  "
$ws.Range("D29").Value = "1.506"
$ws.Range("E29").Value = "  +0.54%  "
$ws.Range("E30").Value = "  -3.52%  "
$ws.Range("D31").Value = "4.401"
$ws.Range("E31").Value = "  +0.29%  "
$ws.Range("D32").Value = "4.265"
$ws.Range("E32").Value = "  +3.70%  "
$ws.Range("D33").Value = "0.05163"
$ws.Range("E33").Value = "  -2.07%  "
$ws.Range("D34").Value = "0.7925"
$ws.Range("E34").Value = "  +8.95%  "
$ws.Range("E35").Value = "  -2.35%  "
$ws.Range("E36").Value = "  -2.84%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01860"
$ws.Range("E38").Value = "  -0.45%  "
$ws.Range("D39").Value = "2.707"
$ws.Range("E39").Value = "  -0.51%  "
$ws.Range("D40").Value = "1.160.42"
$ws.Range("E40").Value = "  -5.61%  "
$ws.Range("D41").Value = "6.319"
$ws.Range("E41").Value = "  +2.76%  "
$ws.Range("D42").Value = "0.8976"
$ws.Range("E42").Value = "  -1.70%  "
$ws.Range("D43").Value = "73.19"
$ws.Range("E43").Value = "  -0.57%  "
$ws.Range("D44").Value = "0.9995"
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").Value = "103.42"
$ws.Range("E45").Value = "  +1.28%  "
$ws.Range("D46").Value = "2.011.47"
$ws.Range("E46").Value = "  -0.32%  "
$ws.Range("D47").Value = "0.5161"
$ws.Range("E47").Value = "  -2.49%  "
$ws.Range("D48").Value = "1.778"
$ws.Range("E48").Value = "  -1.29%  "
$ws.Range("D49").Value = "9.348"
$ws.Range("E49").Value = "  +0.14%  "
$ws.Range("E50").Value = "  -1.19%  "
$ws.Range("D51").Value = "0.4295"
$ws.Range("E51").Value = "  -0.76%  "
